$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-31 Saturday" "2026-02-01 Sunday"

Replace-Text "387÷2=" "810÷5="
Replace-Text "648÷2=" "274÷9="
Replace-Text "456÷5=" "747÷7="
Replace-Text "884÷3=" "186÷7="
Replace-Text "201÷2=" "457÷5="
Replace-Text "985÷5=" "325÷9="
Replace-Text "440÷3=" "192÷5="
Replace-Text "413÷7=" "462÷5="
Replace-Text "136÷5=" "354÷7="
Replace-Text "300÷7=" "921÷7="
Replace-Text "425÷2=" "971÷4="
Replace-Text "273÷6=" "235÷5="
Replace-Text "611÷8=" "622÷2="
Replace-Text "316÷7=" "167÷6="
Replace-Text "532÷3=" "503÷8="
Replace-Text "719÷4=" "176÷7="
Replace-Text "715÷7=" "549÷4="
Replace-Text "954÷5=" "710÷2="
Replace-Text "758÷5=" "741÷5="
Replace-Text "365÷8=" "432÷5="
Replace-Text "720÷2=" "978÷3="
Replace-Text "923÷4=" "357÷8="
Replace-Text "350÷5=" "441÷8="
Replace-Text "306÷8=" "612÷7="
Replace-Text "415÷5=" "542÷6="
